$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 3-8 as part of repulled/recalculated data
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -1
